$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("March 2018 " + [char]0x2013 + " November 2018", $true, $true, $false, $false, $false, `
              $true, 1, $false, "April 2018 " + [char]0x2013 + " November 2018", 2)
